# merge scenarios from master
$wb = $excel.ActiveWorkbook

$wsDescription = $wb.Worksheets.Item("descriptions")

# --- Insert two new rows into the "descriptions" sheet (rows 15 & 16), ---
# --- pushing the existing fishing/length rows down two positions.      ---
$wsDescription.Rows.Item(15).Insert()
$wsDescription.Rows.Item(15).Insert()

# --- Fill in the two brand-new rows (growth + 1yr-forecast scenarios) ---
$wsDescription.Range("A15").Value = "E"
$wsDescription.Range("B15").Value = 990
$wsDescription.Range("C15").Formula = "=CONCATENATE(A15,B15)"
$wsDescription.Range("D15").Value = "fixed - 1 forecast"

$wsDescription.Range("A16").Value = "E"
$wsDescription.Range("B16").Value = 991
$wsDescription.Range("C16").Formula = "=CONCATENATE(A16,B16)"
$wsDescription.Range("D16").Value = "internal - 1 forecast"

# --- Populate the new "long description" column E for every data row ---
# --- (order follows the source merge, not strict row order)           ---
$wsDescription.Range("E2").Value  = "No age comp data for any fleet"
$wsDescription.Range("E3").Value  = "Age comp data for the fishery"
$wsDescription.Range("E4").Value  = "Age comp data for the fishery and every other year of total survey years"
$wsDescription.Range("E5").Value  = "Age comp data for the fishery and for every year of the total survey years"
$wsDescription.Range("E6").Value  = "No conditional age at length data"
$wsDescription.Range("E7").Value  = "Fishery conditional age at length data"
$wsDescription.Range("E9").Value  = "No mean length at age data"
$wsDescription.Range("E10").Value = "Fishery mean length at age data"
$wsDescription.Range("E8").Value  = "Survey conditional age at length data"
$wsDescription.Range("E11").Value = "Survey mean length at age data"
$wsDescription.Range("E12").Value = "all growth parameters fixed at their true values"
$wsDescription.Range("E13").Value = "all growth parameters are estimated internally"
$wsDescription.Range("E14").Value = "all growth parameters are estimated using Christine's external function"
$wsDescription.Range("E15").Value = "all growth parameters are fixed at their true values and one year of forecasting is done"
$wsDescription.Range("E16").Value = "all growth parameters are estimated internally and one year of forecasting is done"
$wsDescription.Range("E17").Value = "Constant fishing at  a percentage of FMSY"
$wsDescription.Range("E18").Value = "Two way trip"
$wsDescription.Range("E19").Value = "A one way trip fishing scenario"
$wsDescription.Range("E21").Value = "Fishery length composition data and every other year of composition data from the survey years"
$wsDescription.Range("E22").Value = "Fishery length composition data and every year of composition data from the survey years"
$wsDescription.Range("E20").Value = "Fishery length composition data"

# --- Make "descriptions" the active sheet (was "scenarios") ---
$wsDescription.Activate()
$null = $wsDescription.Range("E21").Select()
